$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "tricou"
$ws.Range("D1").Value = "Color"
$ws.Range("E1").Value = "Gender"
$ws.Range("B2").Value = 50
$ws.Range("C2").Value = 150
$ws.Range("D2").Value = "alb"
$ws.Range("E2").Value = "barbati"

$ws.Range("A3").Value = "rochie"
$ws.Range("B3").Value = 300
$ws.Range("C3").Value = 700
$ws.Range("E3").Value = "femei"
$ws.Range("D3").Value = "rosu"

$ws.Range("A4").Value = "camasa"
$ws.Range("B4").Value = 170
$ws.Range("C4").Value = 400
$ws.Range("D4").Value = "bleumarin"
$ws.Range("E4").Value = "barbati"

$ws.Range("A5").Value = "geaca"
$ws.Range("B5").Value = 450
$ws.Range("C5").Value = 1200
$ws.Range("D5").Value = "negru"
$ws.Range("E5").Value = "femei"

$ws.Range("A6").Value = "hanorac"
$ws.Range("B6").Value = 250
$ws.Range("C6").Value = 500
$ws.Range("D6").Value = "verde"
$ws.Range("E6").Value = "barbati"

$ws.Range("F6").Select()
